# vocabulary.xlsx edit — fix the module-number boundary between module 11
# and module 12 (the "Spatially lagged variable" row actually belongs to
# module 11, not 12) and highlight the corrected row green.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Module-number corrections (column A) ---------------------------------
# Rows 69-73 were mis-tagged as module 11; they belong to module 12.
$ws.Range("A69").Value = 12
$ws.Range("A70").Value = 12
$ws.Range("A71").Value = 12
$ws.Range("A72").Value = 12
$ws.Range("A73").Value = 12

# Rows 74-79 were mis-tagged as module 12; they belong to module 11.
$ws.Range("A74").Value = 11
$ws.Range("A75").Value = 11
$ws.Range("A76").Value = 11
$ws.Range("A77").Value = 11
$ws.Range("A78").Value = 11
$ws.Range("A79").Value = 11

# --- Formatting -------------------------------------------------------------
# C32 ("Distance" definition) loses its stray formatting override.
$ws.Range("C32").ClearFormats()

# A74 (the row that now correctly starts module 11, "Spatially lagged
# variable") gets highlighted with a solid green fill to flag the fix.
$ws.Range("A74").Interior.Color = 5296274

# --- View / selection state ---------------------------------------------
# Window had scrolled so row 63 was at the top with A49 selected; now it
# scrolls so row 55 is at top and the new last row (A80, one past the used
# range) is selected.
$excel.ActiveWindow.ScrollRow = 55
$ws.Range("A80").Select()
